$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2-18, excluding 6 and 10) were reshuffled: each target row receives
# the Fecha/Calidad/Volumen/Precio minimo/Precio maximo/Precio promedio ponderado/
# Unidad de comercializacion/Precio $/Kg/Kg o Unidades values that used to belong to
# another row, per the permutation below (target row -> source row, both from the
# ORIGINAL worksheet state).
$sourceRows = @{
    2  = 4
    3  = 16
    4  = 5
    5  = 9
    7  = 11
    8  = 12
    9  = 7
    11 = 8
    12 = 15
    13 = 18
    14 = 3
    15 = 17
    16 = 13
    17 = 14
    18 = 2
}

# Capture the original values for every row/column we need before mutating anything,
# since several writes would otherwise clobber data we still need to read later.
$cols = @("D", "I", "J", "K", "L", "M", "N", "P", "Q")
$original = @{}
foreach ($row in $sourceRows.Values) {
    if (-not $original.ContainsKey($row)) {
        $rowData = @{}
        foreach ($col in $cols) {
            $rowData[$col] = $ws.Range("$col$row").Value2
        }
        $original[$row] = $rowData
    }
}

foreach ($target in $sourceRows.Keys) {
    $source = $sourceRows[$target]
    $data = $original[$source]
    foreach ($col in $cols) {
        $ws.Range("$col$target").Value2 = $data[$col]
    }
}
